$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Column C ("Förändrad") holds a date value for rows 2-11.
# Advance that date by one day (45243 -> 45244).
for ($row = 2; $row -le 11; $row++) {
    $cell = $ws.Range("C$row")
    $cell.Value2 = $cell.Value2 + 1
}
